$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1853035143769968
$ws.Range("C2").Value = 0.5686900958466453
$ws.Range("J2").Value = 0.009584664536741214
$ws.Range("P2").Value = 0.1469648562300319
$ws.Range("S2").Value = 0.08945686900958466
$ws.Range("C3").Value = 0.03243243243243243
$ws.Range("J3").Value = 0.01081081081081081
$ws.Range("P3").Value = 0.7621621621621621
$ws.Range("S3").Value = 0.1945945945945946
$ws.Range("J4").Value = 0.06666666666666667
$ws.Range("P4").Value = 0.4444444444444444
$ws.Range("S4").Value = 0.4888888888888889
$ws.Range("B6").Value = 0.05855855855855856
$ws.Range("D6").Value = 0.01801801801801802
$ws.Range("E6").Value = 0.004504504504504504
$ws.Range("F6").Value = 0.04054054054054054
$ws.Range("J6").Value = 0.3243243243243243
$ws.Range("O6").Value = 0.04504504504504504
$ws.Range("Q6").Value = 0.1171171171171171
$ws.Range("R6").Value = 0.07657657657657657
$ws.Range("S6").Value = 0.3153153153153153
$ws.Range("B7").Value = 0.1261682242990654
$ws.Range("D7").Value = 0.03738317757009346
$ws.Range("F7").Value = 0.04672897196261682
$ws.Range("J7").Value = 0.102803738317757
$ws.Range("O7").Value = 0.02336448598130841
$ws.Range("Q7").Value = 0.1588785046728972
$ws.Range("R7").Value = 0.07009345794392523
$ws.Range("S7").Value = 0.4345794392523364
$ws.Range("B8").Value = 0.1256281407035176
$ws.Range("D8").Value = 0.01507537688442211
$ws.Range("F8").Value = 0.06030150753768844
$ws.Range("J8").Value = 0.1231155778894472
$ws.Range("O8").Value = 0.02512562814070352
$ws.Range("Q8").Value = 0.1532663316582915
$ws.Range("R8").Value = 0.09296482412060302
$ws.Range("S8").Value = 0.4045226130653266
$ws.Range("B9").Value = 0.103448275862069
$ws.Range("D9").Value = 0.02758620689655172
$ws.Range("F9").Value = 0.06896551724137931
$ws.Range("J9").Value = 0.103448275862069
$ws.Range("O9").Value = 0.05517241379310345
$ws.Range("Q9").Value = 0.1586206896551724
$ws.Range("R9").Value = 0.06896551724137931
$ws.Range("S9").Value = 0.4137931034482759
$ws.Range("B10").Value = 0.1143790849673203
$ws.Range("D10").Value = 0.0196078431372549
$ws.Range("F10").Value = 0.07516339869281045
$ws.Range("J10").Value = 0.1151960784313725
$ws.Range("O10").Value = 0.02450980392156863
$ws.Range("Q10").Value = 0.2042483660130719
$ws.Range("R10").Value = 0.08169934640522876
$ws.Range("S10").Value = 0.3651960784313725
$ws.Range("F11").Value = 0.003154574132492113
$ws.Range("G11").Value = 0.1545741324921136
$ws.Range("J11").Value = 0.0946372239747634
$ws.Range("K11").Value = 0.2113564668769716
$ws.Range("L11").Value = 0.526813880126183
$ws.Range("S11").Value = 0.009463722397476341
$ws.Range("G12").Value = 0.7660818713450293
$ws.Range("J12").Value = 0.1695906432748538
$ws.Range("L12").Value = 0.02923976608187134
$ws.Range("S12").Value = 0.03508771929824561
$ws.Range("G13").Value = 0.7884615384615384
$ws.Range("J13").Value = 0.2115384615384615
$ws.Range("F15").Value = 0.02145922746781116
$ws.Range("H15").Value = 0.1630901287553648
$ws.Range("I15").Value = 0.03433476394849785
$ws.Range("J15").Value = 0.3218884120171674
$ws.Range("K15").Value = 0.06437768240343347
$ws.Range("N15").Value = 0.004291845493562232
$ws.Range("O15").Value = 0.04721030042918455
$ws.Range("S15").Value = 0.3433476394849785
$ws.Range("F16").Value = 0.01515151515151515
$ws.Range("H16").Value = 0.1464646464646465
$ws.Range("I16").Value = 0.0707070707070707
$ws.Range("J16").Value = 0.4494949494949495
$ws.Range("K16").Value = 0.1262626262626263
$ws.Range("M16").Value = 0.04040404040404041
$ws.Range("O16").Value = 0.04040404040404041
$ws.Range("S16").Value = 0.1111111111111111
$ws.Range("F17").Value = 0.02557544757033248
$ws.Range("H17").Value = 0.1687979539641944
$ws.Range("I17").Value = 0.09207161125319693
$ws.Range("J17").Value = 0.4092071611253197
$ws.Range("K17").Value = 0.09974424552429667
$ws.Range("M17").Value = 0.03324808184143223
$ws.Range("O17").Value = 0.05370843989769821
$ws.Range("S17").Value = 0.1176470588235294
$ws.Range("F18").Value = 0.02906976744186046
$ws.Range("H18").Value = 0.1220930232558139
$ws.Range("I18").Value = 0.06976744186046512
$ws.Range("J18").Value = 0.4651162790697674
$ws.Range("K18").Value = 0.1046511627906977
$ws.Range("M18").Value = 0.02325581395348837
$ws.Range("N18").Value = 0.005813953488372093
$ws.Range("O18").Value = 0.0755813953488372
$ws.Range("S18").Value = 0.1046511627906977
$ws.Range("F19").Value = 0.01976935749588138
$ws.Range("H19").Value = 0.2009884678747941
$ws.Range("I19").Value = 0.06177924217462932
$ws.Range("J19").Value = 0.3764415156507414
$ws.Range("K19").Value = 0.1235584843492586
$ws.Range("M19").Value = 0.02306425041186162
$ws.Range("N19").Value = 0.0008237232289950577
$ws.Range("O19").Value = 0.07907742998352553
$ws.Range("S19").Value = 0.114497528830313
